$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2, 3 and 5 (same market/category/quality) had their
# records cyclically rotated: row2 <- row5, row3 <- row2(old), row5 <- row3(old).
# Row 4 is untouched.

# Row 2 (new values, previously held by row 5)
$ws.Range("D2").Value = 44692
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6

# Row 3 (new values, previously held by row 2)
$ws.Range("D3").Value = 44221
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1300
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1420
$ws.Range("N3").Value = "$/atado"
$ws.Range("O3").Value = "Provincia de Diguillín"
$ws.Range("P3").Value = 1420
$ws.Range("Q3").Value = 1

# Row 5 (new values, previously held by row 3)
$ws.Range("D5").Value = 44691
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3250
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 542
$ws.Range("Q5").Value = 6
